$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ParticipantsTab row's "query" cell (B2) is corrected: a "WITH p" clause is
# inserted right after the WHERE filter (so the later OPTIONAL MATCH chain is
# reseeded only from the participant, not the whole prior row set), and a
# few casing/whitespace nits are cleaned up (WHERE -> where, trailing space
# after RETURN removed, LIMIT -> limit).
$newParticipantsQuery = "MATCH (p:participant)-->(s:study)`n" + `
    "OPTIONAL MATCH (samp:sample)-->(p)`n" + `
    "OPTIONAL MATCH (p)<--(diag:diagnosis)`n" + `
    "OPTIONAL MATCH (samp)<--(f:file)`n" + `
    "OPTIONAL MATCH (f)<--(g:genomic_info)`n" + `
    "WITH s, p, samp, f, g, diag`n" + `
    "where g.library_selection in ['rRNA Depletion']`n" + `
    "WITH p`n" + `
    "OPTIONAL MATCH (p)-->(s:study)`n" + `
    "OPTIONAL MATCH (samp:sample)-->(p)`n" + `
    "WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`n" + `
    "RETURN`n" + `
    "coalesce(p.participant_id,'') as ``Participant ID``,`n" + `
    "coalesce(s.study_name, '') as ``Study Name``,`n" + `
    "coalesce(s.phs_accession,'') as ``Accession``,`n" + `
    "coalesce(p.gender,'') as ``Gender``,`n" + `
    "coalesce(apoc.text.join(samp, ','), '') as ``Samples```n" + `
    "ORDER BY p.participant_id limit 100"

$ws.Range("B2").Value = $newParticipantsQuery

# Reflect the author's new cursor position/selection on the sheet.
$ws.Range("C3").Select()
